# OLX Monitor update — append 8 new listing rows (155-162) to the
# "PODSUMOWANIE" sheet's running log, mirroring rows 147-154 with a
# refreshed check timestamp (2026-02-22 17:34:38) and updated "days seen"
# counts, exactly as produced by the monitoring script's new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# New row data: [targetRow, A timestamp, B profile, C title, D price,
#                E first-seen date, F days-seen, G url, H slug]
$rows = @(
    @(155, "2026-02-22 17:34:38", "poqui", "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza", 2049, "19.12.2025", 65, "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html", "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc"),
    @(156, "2026-02-22 17:34:38", "poqui", "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda", 2299, "19.01.2026", 34, "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html", "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"),
    @(157, "2026-02-22 17:34:38", "poqui", "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy", 2499, "28.10.2025", 117, "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html", "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"),
    @(158, "2026-02-22 17:34:38", "poqui", "Przytulny pokój blisko Politechniki – ul. Przytulna", 549, "10.10.2025", 135, "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html", "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"),
    @(159, "2026-02-22 17:34:38", "pokojewlublinie", "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58", 0, "11.08.2025", 195, "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html", "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"),
    @(160, "2026-02-22 17:34:38", "pokojewlublinie", "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12", 12640, "19.01.2026", 34, "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html", "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"),
    @(161, "2026-02-22 17:34:38", "dawnypatron", "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.", 730, "20.09.2024", 520, "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html", "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"),
    @(162, "2026-02-22 17:34:38", "dawnypatron", "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14", 14690, "05.12.2025", 79, "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html", "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv")
)

foreach ($row in $rows) {
    $targetRow = $row[0]
    $sourceRow = $targetRow - 8

    # Clone the cell formatting (fills/fonts/alignment) from the matching
    # row 8 above, so the new rows keep the same look as the rest of the log.
    $srcRange = $ws.Range("A" + $sourceRow + ":H" + $sourceRow)
    $dstRange = $ws.Range("A" + $targetRow + ":H" + $targetRow)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)

    # Column E holds a plain "DD.MM.YYYY" label, not a real date value, so
    # force it to Text first — otherwise the cell editor reinterprets the
    # string and silently turns it into a date serial number.
    $eCell = $ws.Range("E" + $targetRow)
    $eCell.NumberFormat = "@"

    $ws.Range("A" + $targetRow).Value = $row[1]
    $ws.Range("B" + $targetRow).Value = $row[2]
    $ws.Range("C" + $targetRow).Value = $row[3]
    $ws.Range("D" + $targetRow).Value = $row[4]
    $eCell.Value = $row[5]
    $ws.Range("F" + $targetRow).Value = $row[6]
    $ws.Range("G" + $targetRow).Value = $row[7]
    $ws.Range("H" + $targetRow).Value = $row[8]
}

$excel.CutCopyMode = 0
